# Update bash.lib entry about find usage
#
# The "find a .pdf file" entry (row 6, column C of 工作表1) is expanded
# from a single-line command into a richer note: it now shows both the
# "-name" and the "-type f" forms of `find`, plus a short Traditional
# Chinese explanation and an English postscript about the default
# recursive-search behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Three formatting runs that make up the new cell text.
$run1 = '$find {search-path} -name "*.pdf"' + "`n" + '$find (search-path} -type f              //'
$run2 = "指定查找文件類型文件`n"
$run3 = "p.s. default search into sub-directories."

$full = $run1 + $run2 + $run3

$cell = $ws.Range("C6")
$cell.Value = $full

$len1 = $run1.Length
$len2 = $run2.Length
$len3 = $run3.Length

# Run 1 (the two `find` command lines) keeps the default cell font.

# Run 2: small Chinese explanation, rendered in 細明體 10pt.
$run2Chars = $cell.Characters($len1 + 1, $len2)
$run2Chars.Font.Bold = $false
$run2Chars.Font.Italic = $false
$run2Chars.Font.ColorIndex = -4105
$run2Chars.Font.Name = "細明體"
$run2Chars.Font.Size = 10

# Run 3: English postscript, rendered in Arial 10pt.
$run3Chars = $cell.Characters($len1 + $len2 + 1, $len3)
$run3Chars.Font.Bold = $false
$run3Chars.Font.Italic = $false
$run3Chars.Font.ColorIndex = -4105
$run3Chars.Font.Name = "Arial"
$run3Chars.Font.Size = 10

# Move the active selection to C7, matching the saved view state.
$ws.Range("C7").Select()
